$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "X"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "X"
$ws.Range("R2").Value = ""
$ws.Range("U2").Value = "X"
$ws.Range("W2").Value = "X"
$ws.Range("X2").Value = ""
$ws.Range("D3").Value = "X"
$ws.Range("H3").Value = ""
$ws.Range("K3").Value = ""
$ws.Range("L3").Value = "X"
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("Q3").Value = ""
$ws.Range("T3").Value = ""
$ws.Range("U3").Value = "X"
$ws.Range("W3").Value = "X"
$ws.Range("E4").Value = "X"
$ws.Range("F4").Value = "X"
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = ""
$ws.Range("N4").Value = "X"
$ws.Range("P4").Value = ""
$ws.Range("Q4").Value = ""
$ws.Range("T4").Value = ""
$ws.Range("U4").Value = "X"
$ws.Range("V4").Value = ""
$ws.Range("W4").Value = ""
$ws.Range("X4").Value = ""
$ws.Range("Y4").Value = "X"
$ws.Range("Z4").Value = ""
$ws.Range("AA4").Value = "X"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("N5").Value = ""
$ws.Range("O5").Value = ""
$ws.Range("V5").Value = "X"
$ws.Range("W5").Value = ""
$ws.Range("X5").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = "X"
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = ""
$ws.Range("J6").Value = ""
$ws.Range("L6").Value = ""
$ws.Range("N6").Value = "X"
$ws.Range("P6").Value = ""
$ws.Range("R6").Value = "X"
$ws.Range("V6").Value = "X"
$ws.Range("Z6").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("G7").Value = ""
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = ""
$ws.Range("N7").Value = ""
$ws.Range("O7").Value = ""
$ws.Range("Q7").Value = ""
$ws.Range("R7").Value = "X"
$ws.Range("T7").Value = ""
$ws.Range("X7").Value = ""
$ws.Range("Y7").Value = ""
$ws.Range("Z7").Value = ""
$ws.Range("AB7").Value = "X"
$ws.Range("D8").Value = "X"
$ws.Range("E8").Value = "X"
$ws.Range("F8").Value = ""
$ws.Range("J8").Value = "X"
$ws.Range("K8").Value = ""
$ws.Range("N8").Value = "X"
$ws.Range("T8").Value = ""
$ws.Range("U8").Value = ""
$ws.Range("V8").Value = ""
$ws.Range("X8").Value = "X"
$ws.Range("Y8").Value = ""
$ws.Range("AA8").Value = ""
$ws.Range("AB8").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = "X"
$ws.Range("F9").Value = ""
$ws.Range("G9").Value = ""
$ws.Range("K9").Value = "X"
$ws.Range("M9").Value = ""
$ws.Range("P9").Value = "X"
$ws.Range("Q9").Value = ""
$ws.Range("R9").Value = ""
$ws.Range("U9").Value = ""
$ws.Range("V9").Value = ""
$ws.Range("W9").Value = ""
$ws.Range("X9").Value = ""
$ws.Range("D10").Value = "X"
$ws.Range("F10").Value = ""
$ws.Range("G10").Value = ""
$ws.Range("I10").Value = ""
$ws.Range("J10").Value = "X"
$ws.Range("K10").Value = "X"
$ws.Range("L10").Value = "X"
$ws.Range("P10").Value = "X"
$ws.Range("Q10").Value = ""
$ws.Range("R10").Value = ""
$ws.Range("U10").Value = ""
$ws.Range("X10").Value = "X"
$ws.Range("Z10").Value = ""
$ws.Range("AA10").Value = ""
$ws.Range("AB10").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("G12").Value = ""
$ws.Range("H12").Value = ""
$ws.Range("I12").Value = ""
$ws.Range("K12").Value = ""
$ws.Range("F13").Value = ""
$ws.Range("G13").Value = ""
$ws.Range("H13").Value = ""
$ws.Range("I13").Value = ""
$ws.Range("K13").Value = ""
$ws.Range("I15").Value = ""
$ws.Range("K16").Value = ""
$ws.Range("M16").Value = ""
$ws.Range("V16").Value = ""
$ws.Range("F17").Value = ""
$ws.Range("I17").Value = ""
$ws.Range("M17").Value = ""
$ws.Range("AA17").Value = ""
